$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.301.93"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "1.832.49"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.82%  "
$ws.Range("D5").Value = "'314.40"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("D7").Value = "'0.4734"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("D8").Value = "'0.3686"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").Value = "'0.07443"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").Value = "'0.8858"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").Value = "'20.49"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").Value = "1.886.21"
$ws.Range("E12").Value = "  +3.73%  "
$ws.Range("D13").Value = "'0.07328"
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("D14").Value = "'5.427"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "'94.06"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("D16").Value = "'6.560"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'0.000008792"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").Value = "27.582.98"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").Value = "'14.77"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "'5.286"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").Value = "'10.67"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D24").Value = "2.103.16"
$ws.Range("E24").Value = "  +2.78%  "
$ws.Range("D25").Value = "'1.898"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").Value = "'151.89"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("D28").Value = "'2.137"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "'5.232"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").Value = "'117.32"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").Value = "'0.08990"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").Value = "'0.7492"
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "'4.544"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").Value = "'2.947"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").Value = "'1.095"
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("D38").Value = "'0.05337"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").Value = "'0.01955"
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("D40").Value = "'2.436"
$ws.Range("E40").Value = "  +3.61%  "
$ws.Range("D41").Value = "'2.965"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "'7.246"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "'0.5285"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'0.1658"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "'8.493"
$ws.Range("D46").Value = "'0.4923"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").Value = "'10.52"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'105.14"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("D49").Value = "'1.010"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").Value = "'1.664"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  +0.08%  "
